$wb = $excel.ActiveWorkbook

# Rename sheets (new task-order id suffixes)
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687400167465"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687418612833"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687418642373"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687419349778"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687420115108"

# Sheet 1 - GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687399821029.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687399978545.csv"
$ws1.Range("B4").Value = "go_stims-16511687399999611.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687400147474.csv"

# Sheet 2 - NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16511687411080894.csv"
$ws2.Range("B3").Value = "TB-16511687417600527.csv"
$ws2.Range("B4").Value = "OB-1651168741455329.csv"
$ws2.Range("B5").Value = "TB-16511687418370667.csv"
$ws2.Range("B6").Value = "ZB-match_0-16511687410673926.csv"
$ws2.Range("B7").Value = "TB-165116874156039.csv"
$ws2.Range("B8").Value = "ZB-match_7-16511687404791982.csv"
$ws2.Range("B9").Value = "OB-1651168741251309.csv"
$ws2.Range("B10").Value = "ZB-match_1-16511687409267058.csv"

# Sheet 3 - RS
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 - TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687418889785.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687418672867.csv"
$ws4.Range("B4").Value = "MM_stims-16511687419186788.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687418899803.csv"
$ws4.Range("B6").Value = "MM_stims-16511687419340112.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687419196737.csv"

# Sheet 5 - vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687419392078.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651168741996684.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168741966365.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687419507759.csv"
